# Update crypto price/volume figures on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (D value or $null, E value)
$updates = @(
    @{ Row = 2;  D = "26.789.71";   E = "  +0.55%  " },
    @{ Row = 3;  D = "1.643.93";    E = "  +0.12%  " },
    @{ Row = 4;  D = $null;         E = "  +0.55%  " },
    @{ Row = 5;  D = "216.90";      E = "  +0.83%  " },
    @{ Row = 6;  D = $null;         E = "  -0.47%  " },
    @{ Row = 7;  D = $null;         E = "  +0.48%  " },
    @{ Row = 8;  D = $null;         E = "  -0.07%  " },
    @{ Row = 9;  D = $null;         E = "  +0.21%  " },
    @{ Row = 10; D = "19.17";       E = "  +0.35%  " },
    @{ Row = 11; D = $null;         E = "  +0.45%  " },
    @{ Row = 12; D = "1.661.75";    E = "  +1.43%  " },
    @{ Row = 13; D = $null;         E = "  -0.70%  " },
    @{ Row = 14; D = $null;         E = "  -0.11%  " },
    @{ Row = 15; D = $null;         E = "  -0.74%  " },
    @{ Row = 16; D = "26.777.89";   E = "  +0.27%  " },
    @{ Row = 18; D = "213.92";      E = "  -0.88%  " },
    @{ Row = 19; D = $null;         E = "  +0.48%  " },
    @{ Row = 20; D = $null;         E = "  +0.87%  " },
    @{ Row = 21; D = "2.41";        E = "  +8.91%  " },
    @{ Row = 22; D = "6.24";        E = $null },
    @{ Row = 23; D = "9.31";        E = "  -2.01%  " },
    @{ Row = 24; D = "145.91";      E = "  +0.22%  " },
    @{ Row = 25; D = $null;         E = "  +0.32%  " },
    @{ Row = 26; D = $null;         E = "  -1.59%  " },
    @{ Row = 27; D = "7.18";        E = "  +0.12%  " },
    @{ Row = 28; D = $null;         E = "  -0.43%  " },
    @{ Row = 29; D = $null;         E = "  -1.82%  " },
    @{ Row = 30; D = $null;         E = "  +0.72%  " },
    @{ Row = 31; D = $null;         E = "  -0.81%  " },
    @{ Row = 33; D = "1.283.32";    E = "  +0.22%  " },
    @{ Row = 34; D = $null;         E = "  -0.19%  " },
    @{ Row = 35; D = "2.45";        E = "  +1.46%  " },
    @{ Row = 36; D = $null;         E = "  -0.54%  " },
    @{ Row = 37; D = $null;         E = "  +0.39%  " },
    @{ Row = 38; D = $null;         E = "  -1.55%  " },
    @{ Row = 39; D = $null;         E = "  +0.39%  " },
    @{ Row = 40; D = "0.804";       E = "  -1.74%  " },
    @{ Row = 41; D = $null;         E = "  -0.73%  " },
    @{ Row = 42; D = "5.30";        E = "  -2.75%  " },
    @{ Row = 43; D = "1.783.40";    E = "  +0.02%  " },
    @{ Row = 44; D = "61.36";       E = "  +2.79%  " },
    @{ Row = 45; D = "91.98";       E = "  +1.16%  " },
    @{ Row = 46; D = $null;         E = "  -0.15%  " },
    @{ Row = 47; D = "0.0₆0102";    E = "  -2.23%  " },
    @{ Row = 48; D = $null;         E = "  +0.36%  " },
    @{ Row = 49; D = "7.63";        E = "  -2.16%  " },
    @{ Row = 50; D = $null;         E = "  +0.21%  " },
    @{ Row = 51; D = $null;         E = "  +0.21%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cell = $ws.Range("D$r")
        # Force text interpretation so values like "216.90" or "19.17" are
        # stored verbatim instead of being parsed (and rounded) as numbers.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        # Restore the default "Normal" style so no stray number-format
        # attribute is left behind on the cell.
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Range("E$r").Value = $u.E
    }
}
